$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MenuSheet")

# Update existing rows 9-11 with new menu labels
$ws.Range("A9").Value = "LOUNGEWEAR"
$ws.Range("A10").Value = "MASKS"
$ws.Range("A11").Value = "FAB FREEDOM SALE"

# Add two new rows (12-13), copying the formatting of row 11 (same style as 9-11)
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A12").Value = "GIFT CARD"
$ws.Range("A13").Value = "FAB FIX"

$ws.Range("A13").Select()
